$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8140609860420227
$ws.Range("B1").Value = 1.401088118553162
$ws.Range("D1").Value = 1.775173664093018
$ws.Range("E1").Value = 1.157453656196594
